$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 38464292
$ws.Range("I86").Value = 2160.818
$ws.Range("J86").Value = 66669852
$ws.Range("K86").Value = 2160.818
$ws.Range("L86").Value = 66669852
$ws.Range("M86").Value = -1037.818
$ws.Range("N86").Value = -66672098
$ws.Range("H89").Value = 38464292
$ws.Range("I89").Value = 2160.818
$ws.Range("J89").Value = 66669852
$ws.Range("K89").Value = 10804.09
$ws.Range("L89").Value = 333349260
$ws.Range("M89").Value = -5188.09
$ws.Range("N89").Value = -333360492
$ws.Range("H92").Value = 3718.8
$ws.Range("I92").Value = 3820.889
$ws.Range("K92").Value = 3820.889
$ws.Range("M92").Value = -2572.889
$ws.Range("H137").Value = 1220.5883
$ws.Range("I137").Value = 957.1429
$ws.Range("J137").Value = 2450
$ws.Range("K137").Value = 2871.4287
$ws.Range("L137").Value = 7350
$ws.Range("M137").Value = -321.4287000000004
$ws.Range("N137").Value = -12450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 53958.95
$ws.Range("I2").Value = 84649.5
$ws.Range("J2").Value = 1346.5714
$ws.Range("K2").Value = 84649.5
$ws.Range("L2").Value = 1346.5714
$ws.Range("M2").Value = -84536.5
$ws.Range("N2").Value = -1572.5714
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H97").Value = 432.36
$ws.Range("I97").Value = 387.34784
$ws.Range("K97").Value = 387.34784
$ws.Range("M97").Value = 108.65216
$ws.Range("H116").Value = 53958.95
$ws.Range("I116").Value = 84649.5
$ws.Range("J116").Value = 1346.5714
$ws.Range("K116").Value = 84649.5
$ws.Range("L116").Value = 1346.5714
$ws.Range("M116").Value = -82355.5
$ws.Range("N116").Value = -5934.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 53958.95
$ws.Range("I3").Value = 84649.5
$ws.Range("J3").Value = 1346.5714
$ws.Range("K3").Value = 84649.5
$ws.Range("L3").Value = 1346.5714
$ws.Range("M3").Value = -84535.5
$ws.Range("N3").Value = -1574.5714
$ws.Range("H99").Value = 58825124
$ws.Range("I99").Value = 66668344
$ws.Range("J99").Value = 990
$ws.Range("K99").Value = 66668344
$ws.Range("L99").Value = 990
$ws.Range("M99").Value = -66666846
$ws.Range("N99").Value = -3986
$ws.Range("H134").Value = 24375.979
$ws.Range("I134").Value = 1895.7715
$ws.Range("J134").Value = 103056.7
$ws.Range("K134").Value = 5687.3145
$ws.Range("L134").Value = 309170.1
$ws.Range("M134").Value = -3152.3145
$ws.Range("N134").Value = -314240.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2054.6155
$ws.Range("I132").Value = 1583.5652
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 4750.6956
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -2220.6956
$ws.Range("N132").Value = -22058

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1816.4
$ws.Range("I5").Value = 1680.1111
$ws.Range("K5").Value = 5040.3333
$ws.Range("M5").Value = -4928.3333
$ws.Range("H80").Value = 5075.778
$ws.Range("I80").Value = 1370.5
$ws.Range("J80").Value = 8040
$ws.Range("K80").Value = 4111.5
$ws.Range("L80").Value = 24120
$ws.Range("M80").Value = -3175.5
$ws.Range("N80").Value = -25992
$ws.Range("H83").Value = 5075.778
$ws.Range("I83").Value = 1370.5
$ws.Range("J83").Value = 8040
$ws.Range("K83").Value = 12334.5
$ws.Range("L83").Value = 72360
$ws.Range("M83").Value = -7654.5
$ws.Range("N83").Value = -81720
$ws.Range("H122").Value = 1251693.6
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 1430364.1
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 12873276.9
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -12878176.9
$ws.Range("H132").Value = 1950
$ws.Range("J132").Value = 1950
$ws.Range("L132").Value = 17550
$ws.Range("N132").Value = -22610
$ws.Range("H135").Value = 1816.4
$ws.Range("I135").Value = 1680.1111
$ws.Range("K135").Value = 15120.9999
$ws.Range("M135").Value = -12585.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 9000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H97").Value = 3025.7307
$ws.Range("I97").Value = 3268.0527
$ws.Range("J97").Value = 2368
$ws.Range("K97").Value = 3268.0527
$ws.Range("L97").Value = 2368
$ws.Range("M97").Value = -2772.0527
$ws.Range("N97").Value = -3360
$ws.Range("H132").Value = 3146.1177
$ws.Range("I132").Value = 2660.4614
$ws.Range("J132").Value = 4724.5
$ws.Range("K132").Value = 7981.3842
$ws.Range("L132").Value = 14173.5
$ws.Range("M132").Value = -5451.3842
$ws.Range("N132").Value = -19233.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7884.625
$ws.Range("I61").Value = 8775.286
$ws.Range("J61").Value = 1650
$ws.Range("K61").Value = 8775.286
$ws.Range("L61").Value = 1650
$ws.Range("M61").Value = -8573.286
$ws.Range("N61").Value = -2054
$ws.Range("H100").Value = 10203629
$ws.Range("I100").Value = 14029365
$ws.Range("J100").Value = 1666.6666
$ws.Range("K100").Value = 14029365
$ws.Range("L100").Value = 1666.6666
$ws.Range("M100").Value = -14028824
$ws.Range("N100").Value = -2748.6666
$ws.Range("H113").Value = 7884.625
$ws.Range("I113").Value = 8775.286
$ws.Range("J113").Value = 1650
$ws.Range("K113").Value = 8775.286
$ws.Range("L113").Value = 1650
$ws.Range("M113").Value = -6605.286
$ws.Range("N113").Value = -5990
$ws.Range("H122").Value = 51050
$ws.Range("J122").Value = 2100
$ws.Range("L122").Value = 6300
$ws.Range("N122").Value = -11200
$ws.Range("H132").Value = 1781.2157
$ws.Range("I132").Value = 1580.0834
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 4740.2502
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -2210.2502
$ws.Range("N132").Value = -20058.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2362
$ws.Range("I100").Value = 1548.25
$ws.Range("J100").Value = 2827
$ws.Range("K100").Value = 3096.5
$ws.Range("L100").Value = 5654
$ws.Range("M100").Value = -2555.5
$ws.Range("N100").Value = -6736
$ws.Range("H132").Value = 1615.6154
$ws.Range("I132").Value = 1034.5172
$ws.Range("J132").Value = 3300.8
$ws.Range("K132").Value = 3103.5516
$ws.Range("L132").Value = 9902.400000000001
$ws.Range("M132").Value = -573.5515999999998
$ws.Range("N132").Value = -14962.4

Write-Output "done"